# Applies the "Alvearie" gh-pages deploy edit:
#  - Metadata sheet: Version 5.0.0 -> 6.0.0, Date bump, Publisher filled in,
#    "Contact / No display for ContactDetail" row replaced by a
#    "Jurisdiction / United States of America" row, and the duplicate
#    Contact row removed entirely.
#  - Elements sheet: root Extension row's Short/Definition updated to the
#    resource-specific title/description text.

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")

# Version
$meta.Range("B3").Value = "6.0.0"

# Date
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher
$meta.Range("B9").Value = "Alvearie Team"

# Row 10 was "Contact" / "No display for ContactDetail" -> becomes "Jurisdiction" / "United States of America"
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# Row 11 duplicated the old "Contact" / "No display for ContactDetail" row; remove it so
# everything below shifts up by one (dimension A1:B21 -> A1:B20).
$meta.Rows.Item(11).Delete()

$elements = $wb.Worksheets.Item("Elements")

# Root "Extension" element row: Short/Definition now mirror the StructureDefinition's own
# Title/Description instead of the generic "Extension" / "An Extension" text.
$elements.Range("K2").Value = "Employee Performance Rating"
$elements.Range("L2").Value = "Code indicating the performance rating of the employee, assigned as of the last review period"
